$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "30-jun" column header in T1
$ws.Range("T1").Value = "30-jun"

# Fill in the T2:T18 values for the new column, mirroring column S
$values = @(
    0,
    15.527066483319157,
    13.135707916075548,
    17.126746841624716,
    0,
    7.1497709652943229,
    7.3671837859201563,
    13.628403751588257,
    19.915551216535398,
    10.798594735818151,
    0,
    12.63906448205965,
    0,
    0,
    11.091167731771625,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $values[$i]
}

# Update the selection to reflect the newly added column
$ws.Range("T2:T18").Select()
